$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.497.51'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.913.47'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.38'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4842'
$ws.Range("E7").Value = '  +2.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4071'
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08165'
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("E10").Value = '  +2.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.86'
$ws.Range("E11").Value = '  +5.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.901.84'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.037'
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.193'
$ws.Range("E14").Value = '  +1.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.17'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06779'
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.72'
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.516.56'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.629'
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.181'
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.132.97'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("E26").Value = '  +8.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.46'
$ws.Range("E27").Value = '  +1.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.10'
$ws.Range("E28").Value = '  +1.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.117'
$ws.Range("E29").Value = '  +1.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.77'
$ws.Range("E30").Value = '  +2.60%  '
$ws.Range("E31").Value = '  -3.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09545'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.488'
$ws.Range("E33").Value = '  +2.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.560'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.396'
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02285'
$ws.Range("E36").Value = '  +1.71%  '
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.190'
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("E39").Value = '  +7.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5966'
$ws.Range("E40").Value = '  +2.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.983'
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1856'
$ws.Range("E42").Value = '  +1.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.282'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.399'
$ws.Range("E44").Value = '  -4.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07632'
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.50'
$ws.Range("E46").Value = '  +2.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5586'
$ws.Range("E47").Value = '  +1.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.958'
$ws.Range("E48").Value = '  +2.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.80'
$ws.Range("E49").Value = '  +2.67%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.62'
$ws.Range("E50").Value = '  +2.15%  '
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.412'
$ws.Range("E51").Value = '  +3.09%  '
